$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Login")

# Set new cell D3 to shared string "LMS"
$ws.Range("D3").Value = "LMS"

# Move the active selection to C3 (matches diff's selection change)
$ws.Activate() | Out-Null
$ws.Range("C3").Select() | Out-Null
